# Generate Report for Archive
#
# The localization-status report was regenerated. In the new run, the file
# "86f42771-06db-4da6-93e3-4a9e101966cb" now sorts/lists ahead of
# "fb659db4-70f5-4538-8936-7ccd74a12800" (their rows swap, row 3 <-> row 4)
# across all three sheets (Overview, zh-cn, de-de), and 86f42771's status
# reverts from "Ready for handoff" back to "In Translation" (its own
# handoff file/date are kept). fb659db4 keeps "In Translation" too.
# 089cdb6e (row 2) and e895b03e (row 5) are unaffected.

function Set-HyperlinkDisplay {
    param($ws, $addr, $newText)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $newText
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (zh-cn status), C (de-de
# status), D (Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 3 becomes 86f42771's data
$wsOverview.Range("A3").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.md"
$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"
$wsOverview.Range("D3").Value = "2016-03-22 05:03:26"
Set-HyperlinkDisplay $wsOverview '$A$3' "86f42771-06db-4da6-93e3-4a9e101966cb.md"

# Row 4 becomes fb659db4's data
$wsOverview.Range("A4").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.md"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"
$wsOverview.Range("D4").Value = "2016-03-22 05:02:02"
Set-HyperlinkDisplay $wsOverview '$A$4' "fb659db4-70f5-4538-8936-7ccd74a12800.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A (Source File Name), B (File Extension),
# C (Status), D (Latest Handoff File), E (Latest Handoff Datetime),
# H (Latest Handback DateTime), J (Handoff Reason)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 3 becomes 86f42771's data
$wsZhCn.Range("A3").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("D3").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-22 05:03:18"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J3").Value = "Include"
Set-HyperlinkDisplay $wsZhCn '$A$3' "86f42771-06db-4da6-93e3-4a9e101966cb.md"
Set-HyperlinkDisplay $wsZhCn '$D$3' "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.zh-cn.xlf"

# Row 4 becomes fb659db4's data
$wsZhCn.Range("A4").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "In Translation"
$wsZhCn.Range("D4").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.zh-cn.xlf"
$wsZhCn.Range("E4").Value = "2016-03-22 05:01:47"
$wsZhCn.Range("H4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J4").Value = "Include"
Set-HyperlinkDisplay $wsZhCn '$A$4' "fb659db4-70f5-4538-8936-7ccd74a12800.md"
Set-HyperlinkDisplay $wsZhCn '$D$4' "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de": columns A (Source File Name), B (File Extension),
# C (Status), D (Latest Handoff File), E (Latest Handoff Datetime),
# H (Latest Handback DateTime), J (Handoff Reason)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 3 becomes 86f42771's data
$wsDeDe.Range("A3").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("D3").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-22 05:03:26"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J3").Value = "Include"
Set-HyperlinkDisplay $wsDeDe '$A$3' "86f42771-06db-4da6-93e3-4a9e101966cb.md"
Set-HyperlinkDisplay $wsDeDe '$D$3' "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.de-de.xlf"

# Row 4 becomes fb659db4's data
$wsDeDe.Range("A4").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "In Translation"
$wsDeDe.Range("D4").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.de-de.xlf"
$wsDeDe.Range("E4").Value = "2016-03-22 05:02:02"
$wsDeDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J4").Value = "Include"
Set-HyperlinkDisplay $wsDeDe '$A$4' "fb659db4-70f5-4538-8936-7ccd74a12800.md"
Set-HyperlinkDisplay $wsDeDe '$D$4' "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.de-de.xlf"
